$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 372, shifting existing rows 372:419 down to 373:420
$ws.Rows.Item(372).Insert()

# Populate the newly inserted row 372 with the new record
$ws.Cells.Item(372, 1).Value = 5
$ws.Cells.Item(372, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(372, 3).Value = "Maule"
$ws.Cells.Item(372, 4).Value = 45124
$ws.Cells.Item(372, 4).NumberFormat = $ws.Cells.Item(373, 4).NumberFormat
$ws.Cells.Item(372, 5).Value = 7
$ws.Cells.Item(372, 6).Value = "Fruta"
$ws.Cells.Item(372, 7).Value = 100108
$ws.Cells.Item(372, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(372, 9).Value = 100108005
$ws.Cells.Item(372, 10).Value = "Piña"
$ws.Cells.Item(372, 11).Value = "Caramelo"
$ws.Cells.Item(372, 12).Value = "Segunda"
$ws.Cells.Item(372, 13).Value = 150
$ws.Cells.Item(372, 14).Value = 22000
$ws.Cells.Item(372, 15).Value = 22000
$ws.Cells.Item(372, 16).Value = 22000
$ws.Cells.Item(372, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(372, 18).Value = "Ecuador"
$ws.Cells.Item(372, 19).Value = 1571
$ws.Cells.Item(372, 20).Value = 14
